$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New paragraphs appended below the existing "Khả năng gỡ lỗi..." bullet
# (row 59), mirroring the blank-row spacing pattern already used
# throughout the sheet (one blank separator row, except between the
# intro sentence and the first tool bullet, and between consecutive tool
# bullets, which are contiguous).
$ws.Range("B62").Value = "Dưới đây là một số công cụ tự động hóa hiện đại mà bạn có thể sử dụng thay thế cho tự động hóa trình duyệt IE bằng VBScript và Windows Script:"
$ws.Range("B64").Value = "Selenium WebDriver: Selenium là một công cụ tự động hóa trình duyệt mạnh mẽ và phổ biến. WebDriver cho phép bạn kiểm soát các trình duyệt hiện đại như Chrome, Firefox, Edge, và Safari thông qua mã lập trình trong nhiều ngôn ngữ như Java, Python, C#, và JavaScript. Selenium cung cấp các API mạnh mẽ để tương tác với các phần tử trên trang web, thao tác trình duyệt và thực hiện kiểm thử tự động."
$ws.Range("B65").Value = "Cypress: Cypress là một công cụ kiểm thử tự động dành riêng cho việc kiểm thử trình duyệt. Nó hỗ trợ thực hiện các thao tác tương tác trên trang web và cung cấp các tính năng hiển thị và gỡ lỗi mạnh mẽ."
$ws.Range("B66").Value = "Playwright: Playwright là một thư viện kiểm thử và tự động hóa trình duyệt được phát triển bởi Microsoft. Nó hỗ trợ kiểm thử trình duyệt Chromium, Firefox và WebKit. Playwright có khả năng xử lý các tác vụ phức tạp như kiểm thử đa trình duyệt cùng một lúc và kiểm thử trình duyệt trong các môi trường đa nền tảng."
$ws.Range("B67").Value = "TestCafe: TestCafe là một công cụ tự động hóa kiểm thử trình duyệt hoàn toàn dựa trên JavaScript. Nó cho phép kiểm thử trên các trình duyệt hiện đại như Chrome, Firefox, Edge và Safari."
$ws.Range("B68").Value = "Robot Framework: Robot Framework là một framework kiểm thử và tự động hóa linh hoạt hỗ trợ nhiều ngôn ngữ lập trình. Nó cung cấp các thư viện mở rộng cho việc tự động hóa trình duyệt và kiểm thử ứng dụng."
$ws.Range("B70").Value = "Những công cụ này cung cấp tính năng mạnh mẽ hơn và khả năng tương thích tốt hơn với các trình duyệt hiện đại và các tiêu chuẩn web mới. Trước khi lựa chọn công cụ, hãy cân nhắc yêu cầu cụ thể của dự án và khả năng sử dụng của bạn để đảm bảo rằng công cụ được chọn phù hợp nhất."

# Scroll the view back up (author was near the top of the sheet again
# while reviewing/adding text) and leave the final selection on B35.
$excel.Goto($ws.Range("A3"), $true)
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1

$ws.Range("B35").Select()
